$d = $word.ActiveDocument

# --- locate the "Requisitos" bullet-list paragraph -------------------------
# It is the List Bullet paragraph that immediately follows the paragraph
# whose text is exactly "Requisitos" (a Heading2). Looking this up by
# content (rather than a hard-coded paragraph index) keeps the script
# resilient to any unrelated paragraphs elsewhere in the document.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $cur = $d.Paragraphs.Item($i)
    $curText = $cur.Range.Text.TrimEnd([char]13)
    if ($curText -eq "Requisitos") {
        $targetIndex = $i + 1
    }
}

$p = $d.Paragraphs.Item($targetIndex)

$vt = [char]11   # vertical-tab char == the line break produced by <w:br/>
$cr = [char]13   # paragraph mark character

# --- read the existing requirement lines out of the paragraph --------------
$full = $p.Range.Text.TrimEnd($cr)
$parts = $full.Split($vt)
$items = $parts | Where-Object { $_.Length -gt 0 }

# --- sort the requirement lines alphabetically by their course code --------
# (Sort-Object sorts the whole string, which is equivalent here because
# every line starts with the course code followed by " - ".)
$descSorted = $items | Sort-Object -Descending

$oldStart = $p.Range.Start
$oldEnd = $p.Range.End - 1   # exclude the paragraph mark itself

# --- rebuild the paragraph with the lines in sorted order ------------------
# Inserting each line (in descending order) right at the paragraph's start
# pushes the previous insertions further along, so the final left-to-right
# order ends up ascending. Each InsertBefore call creates its own <w:r> run
# (with its own trailing <w:br/>), matching the original run-per-line shape.
foreach ($s in $descSorted) {
    $insertPos = $p.Range.Start
    $ir = $d.Range($insertPos, $insertPos)
    $ir.InsertBefore($s + $vt)
}

# --- drop the old (now trailing, pre-sort) copies of the lines -------------
$insertedLength = 0
foreach ($s in $items) {
    $insertedLength += $s.Length + 1
}

$delStart = $oldStart + $insertedLength
$delEnd = $oldEnd + $insertedLength
$delRange = $d.Range($delStart, $delEnd)
$delRange.Text = ""

Write-Output "Requisitos list re-sorted."
